$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores numeric-looking values (e.g. "296.84",
# "0.0950") as literal text in the source workbook. Force Text format on
# each Price cell before assigning so Excel does not coerce it to a
# floating-point number (which would lose trailing zeros / exact digits).
# The "Volume(1h)" column (E) already holds values like "  -4.67%  " that
# cannot parse as numbers, so no special handling is needed there.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.969.13"
$ws.Range("E2").Value = "  -4.67%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.221.10"
$ws.Range("E3").Value = "  -6.64%  "

# Row 4
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "296.84"
$ws.Range("E5").Value = "  -5.29%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "80.17"
$ws.Range("E6").Value = "  -9.28%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.504"
$ws.Range("E7").Value = "  -4.93%  "

# Row 8
$ws.Range("E8").Value = "  -0.11%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.458"
$ws.Range("E9").Value = "  -7.47%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0774"
$ws.Range("E10").Value = "  -7.92%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "28.03"
$ws.Range("E11").Value = "  -9.49%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.74"
$ws.Range("E12").Value = "  -13.68%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.108"
$ws.Range("E13").Value = "  -1.65%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.556.43"
$ws.Range("E14").Value = "  -6.90%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.11"
$ws.Range("E15").Value = "  -7.15%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.93"
$ws.Range("E16").Value = "  -7.41%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.246.62"
$ws.Range("E17").Value = "  -5.64%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.712"
$ws.Range("E18").Value = "  -6.90%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "38.855.36"
$ws.Range("E19").Value = "  -4.82%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0858"
$ws.Range("E20").Value = "  -6.02%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.72"
$ws.Range("E21").Value = "  -7.22%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.87"
$ws.Range("E22").Value = "  -6.72%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.81"
$ws.Range("E23").Value = "  -9.32%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "225.52"
$ws.Range("E24").Value = "  -3.36%  "

# Row 25
$ws.Range("E25").Value = "  +0.03%  "

# Row 26
$ws.Range("E26").Value = "  -10.16%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.75"
$ws.Range("E27").Value = "  -3.84%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.30"
$ws.Range("E28").Value = "  -6.44%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.12"
$ws.Range("E29").Value = "  -3.87%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.93"
$ws.Range("E30").Value = "  -4.52%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "147.72"
$ws.Range("E31").Value = "  -3.89%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.14"
$ws.Range("E32").Value = "  -8.25%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.30%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.75"
$ws.Range("E34").Value = "  -8.89%  "

# Row 35
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0687"
$ws.Range("E35").Value = "  -6.22%  "

# Row 36
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.31"
$ws.Range("E36").Value = "  -5.48%  "

# Row 37
$ws.Range("E37").Value = "  -4.97%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.65"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0950"
$ws.Range("E39").Value = "  -4.93%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "14.53"
$ws.Range("E40").Value = "  -8.88%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.59"
$ws.Range("E41").Value = "  -7.68%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.59"
$ws.Range("E42").Value = "  -6.46%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.908.97"
$ws.Range("E43").Value = "  -2.50%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.07"
$ws.Range("E44").Value = "  -11.98%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0254"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.14"
$ws.Range("E46").Value = "  -8.72%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.86"
$ws.Range("E47").Value = "  -6.24%  "

# Row 48
$ws.Range("E48").Value = "  -9.37%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.421.71"
$ws.Range("E49").Value = "  -7.03%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "70.73"
$ws.Range("E50").Value = "  -2.99%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "87.12"
$ws.Range("E51").Value = "  -7.10%  "
